$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.668.11"
$ws.Range("E2").Value = "  -1.67%  "
$ws.Range("D3").Value = "2.631.90"
$ws.Range("E3").Value = "  -0.35%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "520.01"
$ws.Range("E5").Value = "  +0.39%  "
$ws.Range("D6").Value = "144.97"
$ws.Range("E6").Value = "  -3.80%  "
$ws.Range("E7").Value = "  +0.38%  "
$ws.Range("D8").Value = "0.575"
$ws.Range("E8").Value = "  -0.43%  "
$ws.Range("D9").Value = "2.642.57"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("E10").Value = "  -3.98%  "
$ws.Range("E11").Value = "  -2.40%  "
$ws.Range("D12").Value = "0.333"
$ws.Range("E12").Value = "  -2.78%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "3.094.57"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "58.698.91"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "20.76"
$ws.Range("E16").Value = "  -3.61%  "
$ws.Range("E17").Value = "  -3.03%  "
$ws.Range("D18").Value = "2.637.38"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").Value = "347.97"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "4.44"
$ws.Range("E20").Value = "  -4.51%  "
$ws.Range("E21").Value = "  -4.10%  "
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "61.65"
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("D25").Value = "0.413"
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("D26").Value = "0.163"
$ws.Range("E26").Value = "  +0.50%  "
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  +0.32%  "
$ws.Range("D28").Value = "0.0₃0803"
$ws.Range("E28").Value = "  -5.03%  "
$ws.Range("E29").Value = "  -2.81%  "
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.13%  "
$ws.Range("E31").Value = "  -5.72%  "
$ws.Range("D32").Value = "18.84"
$ws.Range("E32").Value = "  -1.80%  "
$ws.Range("E33").Value = "  -0.29%  "
$ws.Range("D34").Value = "148.92"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").Value = "0.986"
$ws.Range("E35").Value = "  -6.07%  "
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").Value = "36.61"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("D39").Value = "0.838"
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("D41").Value = "3.61"
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("D42").Value = "279.21"
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "0.0983"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("D45").Value = "0.600"
$ws.Range("E45").Value = "  -4.82%  "
$ws.Range("D46").Value = "19.58"
$ws.Range("E46").Value = "  -1.16%  "
$ws.Range("E47").Value = "  -5.08%  "
$ws.Range("D48").Value = "10.30"
$ws.Range("E48").Value = "  +0.40%  "
$ws.Range("D49").Value = "0.0228"
$ws.Range("E49").Value = "  -2.61%  "
$ws.Range("D50").Value = "1.988.45"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "4.63"
$ws.Range("E51").Value = "  -4.16%  "
